# Auto-applies the "Actualizacion automatica" monthly rollover edit
# described by the target diff:
#  - Sheet "VENTAS POR GRUPO": zero out this period's per-category sales
#    figures (rows 2-55) and reset the row 56 "x de 54" counters to 0.
#  - Sheet "VENTA MENSUAL": roll the monthly columns left by one month
#    (sep->oct->nov->dic->ene headers), shifting each client's C/D/E
#    values accordingly, posting the new "diciembre" actuals, zeroing
#    the new (empty) "enero" column, and adjusting the affected column
#    widths to match the new header text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" -- zero out this period's category sales
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCells = @(
    "M3", "O3", "C4", "D4", "E4", "M4", "N4", "D5", "L5", "M5", "D7", "K7", "M7", "I10", "K10",
    "C12", "K12", "L12", "M13", "M16", "M18", "I19", "L19", "M19", "E24", "M24", "H26", "M26",
    "M28", "D29", "E29", "M29", "O29", "I36", "K36", "M36", "O36", "P36", "C37", "I37", "L37",
    "M37", "M40", "D48", "I48", "D50", "D52", "O52", "H53", "I53", "M53"
)
foreach ($cellRef in $zeroCells) {
    $wsGrupo.Range($cellRef).Value = 0
}

# Row 56 holds "<n> de 54" coverage counters per category; every category
# that lost its only non-zero contributor above now reads "0 de 54".
$row56Cols = @("C", "D", "E", "H", "I", "K", "L", "M", "N", "O", "P")
foreach ($col in $row56Cols) {
    $wsGrupo.Range($col + "56").Value = "0 de 54"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" -- roll the 4-month window forward by one
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Header row: septiembre/octubre/noviembre/diciembre -> octubre/noviembre/diciembre/enero
$wsMensual.Range("C1").Value = "octubre"
$wsMensual.Range("D1").Value = "noviembre"
$wsMensual.Range("E1").Value = "diciembre"
$wsMensual.Range("F1").Value = "enero"

# Column widths follow the new header labels (auto-fit deltas baked into the template).
$wsMensual.Columns.Item(3).ColumnWidth = 13.166666666666666  # "octubre"   -> width 14
$wsMensual.Columns.Item(4).ColumnWidth = 14.166666666666666  # "noviembre" -> width 15
$wsMensual.Columns.Item(6).ColumnWidth = 10.166666666666666  # "enero"     -> width 11

# Per-client monthly figures: shift C<-D<-E, post the new diciembre actual
# into E, and zero the brand-new enero column (F).
$monthlyData = @{
    3 = @{ "C"="1855.16"; "D"="498.96"; "E"="1237.68"; "F"="0" }
    4 = @{ "C"="9124.43"; "D"="12438.61"; "E"="13723.92"; "F"="0" }
    5 = @{ "C"="13542.31"; "D"="3386.11"; "E"="10074.14"; "F"="0" }
    6 = @{ "C"="0" }
    7 = @{ "C"="1453.68"; "D"="1160.13"; "E"="4094.47"; "F"="0" }
    8 = @{ "C"="0"; "D"="314.8"; "E"="305.28" }
    9 = @{ "C"="3172.12"; "D"="924.9400000000001"; "E"="0" }
    10 = @{ "C"="1726.02"; "D"="0"; "E"="560.16"; "F"="0" }
    11 = @{ "C"="2370.29"; "D"="0"; "E"="1324.8" }
    12 = @{ "C"="475.2"; "D"="4643.83"; "E"="3304.64"; "F"="0" }
    13 = @{ "C"="7536.8"; "D"="4486.61"; "E"="3046.3"; "F"="0" }
    15 = @{ "C"="812.16"; "D"="0"; "E"="1417.42" }
    16 = @{ "C"="6465.16"; "D"="21181.89"; "E"="6914.1"; "F"="0" }
    17 = @{ "C"="-354.36"; "D"="1128.58"; "E"="0" }
    18 = @{ "C"="0"; "D"="3299.06"; "E"="2505.46"; "F"="0" }
    19 = @{ "C"="0"; "E"="3219.8"; "F"="0" }
    22 = @{ "D"="140.76"; "E"="0" }
    23 = @{ "C"="0" }
    24 = @{ "C"="7750.68"; "D"="5478.92"; "E"="12621.68"; "F"="0" }
    26 = @{ "C"="0"; "E"="1457.67"; "F"="0" }
    28 = @{ "C"="0"; "D"="5224.76"; "E"="2378.65"; "F"="0" }
    29 = @{ "C"="12734.33"; "D"="8818.120000000001"; "E"="8513.129999999999"; "F"="0" }
    31 = @{ "D"="79.2"; "E"="0" }
    32 = @{ "C"="-21.6"; "D"="7036.98"; "E"="270.92" }
    36 = @{ "C"="17640.85"; "D"="9092.17"; "E"="17009.67"; "F"="0" }
    37 = @{ "C"="8317.59"; "D"="3690.44"; "E"="7074.68"; "F"="0" }
    38 = @{ "C"="258.08"; "D"="0" }
    39 = @{ "D"="121.31"; "E"="0" }
    40 = @{ "C"="9434.389999999999"; "D"="0"; "E"="357.5"; "F"="0" }
    41 = @{ "C"="45.69"; "D"="641.5"; "E"="32.53" }
    42 = @{ "C"="929.16"; "D"="0" }
    46 = @{ "C"="457.92"; "D"="0" }
    47 = @{ "C"="0"; "D"="1910.7"; "E"="794.5599999999999" }
    48 = @{ "C"="81"; "D"="2139.68"; "E"="1213.49"; "F"="0" }
    49 = @{ "C"="189.19"; "D"="0" }
    50 = @{ "C"="55.65"; "D"="0"; "E"="91.58"; "F"="0" }
    53 = @{ "C"="5829.37"; "D"="2479.21"; "E"="3511.8"; "F"="0" }
    54 = @{ "C"="5829.37"; "D"="2479.21"; "E"="3511.8"; "F"="0" }
    55 = @{ "C"="264.06"; "D"="1011.96"; "E"="1657.97"; "F"="0" }
    56 = @{ "C"="264.06"; "D"="1011.96"; "E"="1657.97"; "F"="0" }
    60 = @{ "C"="118238.76"; "D"="104820.4"; "E"="113883.77"; "F"="0" }
}

foreach ($rowNum in $monthlyData.Keys) {
    $rowVals = $monthlyData[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $wsMensual.Range($col + $rowNum).Value = [double]$rowVals[$col]
    }
}

